# Insert a new data row at row 19 (shifts existing rows 19-79 down to 20-80)
# and populate it with the new "Poroto granado" price record for
# Terminal Hortofrutícola Agro Chillán, matching the committed diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(19).Insert()

$ws.Cells.Item(19, 1).Value  = 7
$ws.Cells.Item(19, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(19, 3).Value  = "Ñuble"
$ws.Cells.Item(19, 4).Value  = 44592
$ws.Cells.Item(19, 5).Value  = 16
$ws.Cells.Item(19, 6).Value  = 100112030
$ws.Cells.Item(19, 7).Value  = "Poroto granado"
$ws.Cells.Item(19, 8).Value  = "Sin especificar"
$ws.Cells.Item(19, 9).Value  = "Primera"
$ws.Cells.Item(19, 10).Value = 100
$ws.Cells.Item(19, 11).Value = 23000
$ws.Cells.Item(19, 12).Value = 24000
$ws.Cells.Item(19, 13).Value = 23500
$ws.Cells.Item(19, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(19, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(19, 16).Value = 940
$ws.Cells.Item(19, 17).Value = 25
$ws.Cells.Item(19, 18).Value = "Hortaliza"
